$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25:308 down to 26:309
$ws.Rows(25).Insert()

# Populate the newly inserted row 25 with a new data record (same shape as the
# other rows in this table), with new Fecha (D) and Volumen (J) values.
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = "Femacal de La Calera"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44685
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 100112039
$ws.Range("G25").Value = "Ciboulette"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = 1500
$ws.Range("N25").Value = "`$/docena de atados"
$ws.Range("O25").Value = "Provincia de Quillota"
$ws.Range("P25").Value = 500
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = "Hortaliza"
